$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 78214
$ws.Range("B3").Value = 27.53176515713299
$ws.Range("B4").Value = 2.388879773628455
$ws.Range("B9").Value = 48.42
